$d = $word.ActiveDocument

$d.Content.Find.Execute("82÷7=11, 5", $true, $false, $false, $false, $false, $true, 1, $false, "68÷7=9, 5", 2) | Out-Null
$d.Content.Find.Execute("51÷2=25, 1", $true, $false, $false, $false, $false, $true, 1, $false, "19÷8=2, 3", 2) | Out-Null
$d.Content.Find.Execute("16÷4=4, 0", $true, $false, $false, $false, $false, $true, 1, $false, "67÷5=13, 2", 2) | Out-Null
$d.Content.Find.Execute("49÷9=5, 4", $true, $false, $false, $false, $false, $true, 1, $false, "71÷7=10, 1", 2) | Out-Null
$d.Content.Find.Execute("86÷8=10, 6", $true, $false, $false, $false, $false, $true, 1, $false, "14÷7=2, 0", 2) | Out-Null
$d.Content.Find.Execute("32÷7=4, 4", $true, $false, $false, $false, $false, $true, 1, $false, "64÷3=21, 1", 2) | Out-Null
$d.Content.Find.Execute("25÷4=6, 1", $true, $false, $false, $false, $false, $true, 1, $false, "30÷9=3, 3", 2) | Out-Null
$d.Content.Find.Execute("63÷7=9, 0", $true, $false, $false, $false, $false, $true, 1, $false, "26÷6=4, 2", 2) | Out-Null
$d.Content.Find.Execute("74÷6=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "94÷8=11, 6", 2) | Out-Null
$d.Content.Find.Execute("61÷2=30, 1", $true, $false, $false, $false, $false, $true, 1, $false, "82÷5=16, 2", 2) | Out-Null
$d.Content.Find.Execute("61÷4=15, 1", $true, $false, $false, $false, $false, $true, 1, $false, "36÷7=5, 1", 2) | Out-Null
$d.Content.Find.Execute("98÷8=12, 2", $true, $false, $false, $false, $false, $true, 1, $false, "53÷2=26, 1", 2) | Out-Null
$d.Content.Find.Execute("56÷5=11, 1", $true, $false, $false, $false, $false, $true, 1, $false, "47÷2=23, 1", 2) | Out-Null
$d.Content.Find.Execute("48÷8=6, 0", $true, $false, $false, $false, $false, $true, 1, $false, "39÷9=4, 3", 2) | Out-Null
$d.Content.Find.Execute("33÷7=4, 5", $true, $false, $false, $false, $false, $true, 1, $false, "10÷3=3, 1", 2) | Out-Null
$d.Content.Find.Execute("48÷5=9, 3", $true, $false, $false, $false, $false, $true, 1, $false, "31÷2=15, 1", 2) | Out-Null
$d.Content.Find.Execute("66÷2=33, 0", $true, $false, $false, $false, $false, $true, 1, $false, "93÷6=15, 3", 2) | Out-Null
$d.Content.Find.Execute("86÷6=14, 2", $true, $false, $false, $false, $false, $true, 1, $false, "28÷9=3, 1", 2) | Out-Null
$d.Content.Find.Execute("31÷3=10, 1", $true, $false, $false, $false, $false, $true, 1, $false, "65÷9=7, 2", 2) | Out-Null
$d.Content.Find.Execute("92÷5=18, 2", $true, $false, $false, $false, $false, $true, 1, $false, "13÷4=3, 1", 2) | Out-Null
$d.Content.Find.Execute("67÷8=8, 3", $true, $false, $false, $false, $false, $true, 1, $false, "50÷3=16, 2", 2) | Out-Null
$d.Content.Find.Execute("38÷7=5, 3", $true, $false, $false, $false, $false, $true, 1, $false, "39÷9=4, 3", 2) | Out-Null
$d.Content.Find.Execute("79÷8=9, 7", $true, $false, $false, $false, $false, $true, 1, $false, "54÷4=13, 2", 2) | Out-Null
$d.Content.Find.Execute("74÷4=18, 2", $true, $false, $false, $false, $false, $true, 1, $false, "17÷3=5, 2", 2) | Out-Null
$d.Content.Find.Execute("99÷8=12, 3", $true, $false, $false, $false, $false, $true, 1, $false, "34÷3=11, 1", 2) | Out-Null
